$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of data (row 46 and 47), filling column by column
# so that new shared-string entries are created in the same order as the
# original authoring (A46, A47, B46, B47, C46, C47).
$ws.Range("A46").Value = "45.jpg"
$ws.Range("A47").Value = "46.jpg"

$ws.Range("B46").Value = 1
$ws.Range("B47").Value = 1

$ws.Range("C46").Value = "Зажигалка в нож"
$ws.Range("C47").Value = "Зажигалка в транспортир"

# Update the selected cell to D1
$ws.Range("D1").Select()
